# multi browser implementation for chrome and firefox and accounting cash
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Summary": update a couple of values and drop the now-unused
# "placeholder" rows 7-10, pulling the trailing footer rows up.
# ------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("F2").Value = 1581.47
$wsSummary.Range("A3").Value = 685.07
$wsSummary.Range("E3").Value = 685.07
$wsSummary.Range("F3").Value = 193.97

# Remove rows 6 through 10 (row 6 is already blank) so the old rows
# 11/12 shift up to become rows 6/7.
$wsSummary.Range("A6:A10").EntireRow.Delete()

[void]$wsSummary.Range("E13").Select()

# ------------------------------------------------------------------
# Sheet "Repayment Schedule": refresh the recalculated repayment
# figures and move the active selection.
# ------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

$wsSchedule.Range("F3").Value = 785.8
$wsSchedule.Range("G3").Value = 9214.2000000000007
$wsSchedule.Range("H3").Value = 101.92
$wsSchedule.Range("K3").Value = 887.72
$wsSchedule.Range("P3").Value = 887.72

$wsSchedule.Range("F4").Value = 795.67
$wsSchedule.Range("G4").Value = 8418.5300000000007
$wsSchedule.Range("H4").Value = 92.05
$wsSchedule.Range("K4").Value = 887.72
$wsSchedule.Range("P4").Value = 887.72

$wsSchedule.Range("F5").Value = 785.8
$wsSchedule.Range("G5").Value = 7632.73
$wsSchedule.Range("H5").Value = 101.92
$wsSchedule.Range("K5").Value = 887.72
$wsSchedule.Range("P5").Value = 887.72

$wsSchedule.Range("F6").Value = 812.44
$wsSchedule.Range("G6").Value = 6820.29
$wsSchedule.Range("H6").Value = 75.28
$wsSchedule.Range("K6").Value = 887.72
$wsSchedule.Range("P6").Value = 887.72

$wsSchedule.Range("F7").Value = 818.21
$wsSchedule.Range("G7").Value = 6002.08
$wsSchedule.Range("H7").Value = 69.510000000000005
$wsSchedule.Range("K7").Value = 887.72
$wsSchedule.Range("P7").Value = 887.72

$wsSchedule.Range("F8").Value = 828.52
$wsSchedule.Range("G8").Value = 5173.5600000000004
$wsSchedule.Range("H8").Value = 59.2
$wsSchedule.Range("K8").Value = 887.72
$wsSchedule.Range("P8").Value = 887.72

$wsSchedule.Range("F9").Value = 834.99
$wsSchedule.Range("G9").Value = 4338.57
$wsSchedule.Range("H9").Value = 52.73
$wsSchedule.Range("K9").Value = 887.72
$wsSchedule.Range("P9").Value = 887.72

$wsSchedule.Range("F10").Value = 843.5
$wsSchedule.Range("G10").Value = 3495.07
$wsSchedule.Range("H10").Value = 44.22
$wsSchedule.Range("K10").Value = 887.72
$wsSchedule.Range("P10").Value = 887.72

$wsSchedule.Range("F11").Value = 853.25
$wsSchedule.Range("G11").Value = 2641.82
$wsSchedule.Range("H11").Value = 34.47
$wsSchedule.Range("K11").Value = 887.72
$wsSchedule.Range("P11").Value = 887.72

$wsSchedule.Range("F12").Value = 860.8
$wsSchedule.Range("G12").Value = 1781.02
$wsSchedule.Range("H12").Value = 26.92
$wsSchedule.Range("K12").Value = 887.72
$wsSchedule.Range("P12").Value = 887.72

$wsSchedule.Range("F13").Value = 870.15
$wsSchedule.Range("G13").Value = 910.87
$wsSchedule.Range("H13").Value = 17.57
$wsSchedule.Range("K13").Value = 887.72
$wsSchedule.Range("P13").Value = 887.72

$wsSchedule.Range("F14").Value = 910.87
$wsSchedule.Range("H14").Value = 9.2799999999999994
$wsSchedule.Range("K14").Value = 920.15
$wsSchedule.Range("P14").Value = 920.15

[void]$wsSchedule.Range("E20").Select()

Write-Host "applied repayment schedule + summary updates"
